$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: 2020.05.08 / look at CICM UI to prep for Konrad call, think about UI ---
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "2020.05.08"
$ws.Range("A21").Style = "Normal"

$ws.Range("B21").Value = 0.43055555555555558
$ws.Range("B21").NumberFormat = "h:mm"

$ws.Range("C21").Value = 0.51388888888888895
$ws.Range("C21").NumberFormat = "h:mm"

$ws.Range("D21").Formula = "=C21-B21"
$ws.Range("D21").NumberFormat = "h:mm"

$ws.Range("E21").WrapText = $true
$ws.Range("E21").Value = "look at CICM UI to prep for Konrad call, think about UI"

# --- Row 22: call with Konrad ---
$ws.Range("B22").Value = 0.64583333333333337
$ws.Range("B22").NumberFormat = "h:mm"

$ws.Range("C22").Value = 0.69791666666666663
$ws.Range("C22").NumberFormat = "h:mm"

$ws.Range("D22").Formula = "=C22-B22"
$ws.Range("D22").NumberFormat = "h:mm"

$ws.Range("E22").WrapText = $true
$ws.Range("E22").Value = "call with Konrad"

# --- Row 23: continue working on porting conversion functions ---
$ws.Range("B23").Value = 0.71875
$ws.Range("B23").NumberFormat = "h:mm"

$ws.Range("C23").Value = 0.80208333333333337
$ws.Range("C23").NumberFormat = "h:mm"

$ws.Range("D23").Formula = "=C23-B23"
$ws.Range("D23").NumberFormat = "h:mm"

$ws.Range("E23").WrapText = $true
$ws.Range("E23").Value = "continue working on porting conversion functions"

# --- Row 24 (taller row, ht=30): linker error explanation ---
$ws.Rows.Item(24).RowHeight = 30

$ws.Range("B24").Value = 0.84375
$ws.Range("B24").NumberFormat = "h:mm"

$ws.Range("C24").Value = 0.91666666666666663
$ws.Range("C24").NumberFormat = "h:mm"

$ws.Range("D24").Formula = "=C24-B24"
$ws.Range("D24").NumberFormat = "h:mm"

$ws.Range("E24").WrapText = $true
$ws.Range("E24").Value = "linker error with conversions. Issue was that all the min stuff doesn't use implementation files, so the Only Once rule being violated in the resulting object files. "

# --- Row 25: extend formula (E25 already blank/selected) ---
$ws.Range("D25").Formula = "=C25-B25"
$ws.Range("D25").NumberFormat = "h:mm"

# --- Row 26 (new row): extend formula further ---
$ws.Range("D26").Formula = "=C26-B26"
$ws.Range("D26").NumberFormat = "h:mm"

# --- Scroll the view toward the top of the log (best effort; engine may not persist topLeftCell) ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E25").Select()
